# Updates the crypto price/volume table to reflect the latest scrape.
# Row 30/31 (Kaspa <-> OKB) swapped position in the source ranking, so
# the whole row (Coin, Link, Price, Volume(1h)) is rewritten for those two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.931.70"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3
$ws.Range("D3").Value = "2.841.15"
$ws.Range("E3").Value = "  +2.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'361.26"
$ws.Range("E5").Value = "  +6.52%  "

# Row 6
$ws.Range("D6").Value = "'113.34"
$ws.Range("E6").Value = "  -2.78%  "

# Row 7
$ws.Range("D7").Value = "'0.566"
$ws.Range("E7").Value = "  +4.68%  "

# Row 8
$ws.Range("E8").Value = "  +0.11%  "

# Row 9
$ws.Range("D9").Value = "'0.604"
$ws.Range("E9").Value = "  +4.90%  "

# Row 10
$ws.Range("D10").Value = "'41.67"
$ws.Range("E10").Value = "  -1.43%  "

# Row 11
$ws.Range("D11").Value = "'0.0858"
$ws.Range("E11").Value = "  -1.10%  "

# Row 12
$ws.Range("E12").Value = "  +1.22%  "

# Row 13
$ws.Range("D13").Value = "'19.89"
$ws.Range("E13").Value = "  -0.85%  "

# Row 14
$ws.Range("D14").Value = "'7.78"
$ws.Range("E14").Value = "  +1.92%  "

# Row 15
$ws.Range("D15").Value = "3.290.39"
$ws.Range("E15").Value = "  +2.26%  "

# Row 16
$ws.Range("D16").Value = "2.849.31"
$ws.Range("E16").Value = "  +2.30%  "

# Row 17
$ws.Range("D17").Value = "'0.902"
$ws.Range("E17").Value = "  +1.81%  "

# Row 18
$ws.Range("D18").Value = "51.890.02"
$ws.Range("E18").Value = "  +0.29%  "

# Row 19
$ws.Range("D19").Value = "'7.46"
$ws.Range("E19").Value = "  +7.35%  "

# Row 20
$ws.Range("D20").Value = "'3.17"
$ws.Range("E20").Value = "  -2.78%  "

# Row 21
$ws.Range("D21").Value = "'13.54"
$ws.Range("E21").Value = "  +0.80%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0991"
$ws.Range("E22").Value = "  +1.36%  "

# Row 23
$ws.Range("D23").Value = "'69.91"
$ws.Range("E23").Value = "  -0.33%  "

# Row 24
$ws.Range("D24").Value = "'267.00"
$ws.Range("E24").Value = "  -3.78%  "

# Row 25
$ws.Range("D25").Value = "'2.85"
$ws.Range("E25").Value = "  +3.11%  "

# Row 26
$ws.Range("D26").Value = "'27.10"
$ws.Range("E26").Value = "  +1.07%  "

# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").Value = "'10.40"
$ws.Range("E28").Value = "  +2.03%  "

# Row 30
$ws.Range("B30").Value = "OKB"
$ws.Range("C30").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D30").Value = "'53.48"
$ws.Range("E30").Value = "  +6.64%  "

# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.140"
$ws.Range("E31").Value = "  -1.13%  "

# Row 32
$ws.Range("D32").Value = "'33.89"
$ws.Range("E32").Value = "  -3.47%  "

# Row 33
$ws.Range("D33").Value = "'5.90"
$ws.Range("E33").Value = "  +4.41%  "

# Row 34
$ws.Range("D34").Value = "'0.0445"
$ws.Range("E34").Value = "  +21.97%  "

# Row 35
$ws.Range("D35").Value = "'0.0838"
$ws.Range("E35").Value = "  +2.06%  "

# Row 36
$ws.Range("D36").Value = "'5.21"
$ws.Range("E36").Value = "  +4.71%  "

# Row 37
$ws.Range("E37").Value = "  -0.01%  "

# Row 38
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "  +1.15%  "

# Row 39
$ws.Range("D39").Value = "'2.08"
$ws.Range("E39").Value = "  -2.32%  "

# Row 40
$ws.Range("D40").Value = "'18.29"
$ws.Range("E40").Value = "  -4.00%  "

# Row 41
$ws.Range("D41").Value = "'24.03"
$ws.Range("E41").Value = "  +2.70%  "

# Row 42
$ws.Range("D42").Value = "'2.57"
$ws.Range("E42").Value = "  -4.69%  "

# Row 43
$ws.Range("D43").Value = "'0.117"
$ws.Range("E43").Value = "  +2.12%  "

# Row 44
$ws.Range("D44").Value = "'127.20"
$ws.Range("E44").Value = "  -0.11%  "

# Row 45
$ws.Range("E45").Value = "  -3.22%  "

# Row 46
$ws.Range("D46").Value = "2.112.85"
$ws.Range("E46").Value = "  +0.50%  "

# Row 47
$ws.Range("D47").Value = "'3.37"
$ws.Range("E47").Value = "  +1.48%  "

# Row 48
$ws.Range("D48").Value = "'2.25"
$ws.Range("E48").Value = "  +0.97%  "

# Row 49
$ws.Range("D49").Value = "'1.01"
$ws.Range("E49").Value = "  +11.25%  "

# Row 50
$ws.Range("D50").Value = "'5.83"
$ws.Range("E50").Value = "  +5.14%  "

# Row 51
$ws.Range("D51").Value = "'9.01"
$ws.Range("E51").Value = "  +1.29%  "
